# Update "想去人数" (F) and "最低票价" (G) figures on the 展览 (Exhibition)
# and 全部类型 (All types) sheets to the newly scraped values.
#
# 展览 sheet: row -> [F (new), G (new or $null to leave unchanged)]
# 全部类型 sheet: same idea, but row numbers are shifted for some rows
# because that sheet has one extra row (a 演出/concert entry) inserted
# around row 20 relative to 展览.

$wb = $excel.ActiveWorkbook

function Update-Row {
    param(
        $ws,
        [int]$row,
        $newF,
        $newG
    )

    if ($null -ne $newF) {
        $ws.Cells.Item($row, 6).Value = $newF
    }
    if ($null -ne $newG) {
        $ws.Cells.Item($row, 7).Value = $newG
    }
}

# ---------------------------------------------------------------------
# Sheet "展览" (index 1)
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

Update-Row $wsExpo 2  15211 $null
Update-Row $wsExpo 3  19641 105
Update-Row $wsExpo 5  190   $null
Update-Row $wsExpo 14 232   "已售罄"
Update-Row $wsExpo 15 261   $null
Update-Row $wsExpo 16 79    $null
Update-Row $wsExpo 17 1543  $null
Update-Row $wsExpo 20 123   $null
Update-Row $wsExpo 21 254   $null
Update-Row $wsExpo 22 8309  $null
Update-Row $wsExpo 24 48    $null
Update-Row $wsExpo 27 1285  $null
Update-Row $wsExpo 28 40    $null
Update-Row $wsExpo 31 6696  $null
Update-Row $wsExpo 33 82    $null
Update-Row $wsExpo 34 195   $null
Update-Row $wsExpo 37 5703  $null
Update-Row $wsExpo 41 68    $null

# ---------------------------------------------------------------------
# Sheet "全部类型" (index 4)
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

Update-Row $wsAll 2  15211 $null
Update-Row $wsAll 3  19641 105
Update-Row $wsAll 5  190   $null
Update-Row $wsAll 14 232   "已售罄"
Update-Row $wsAll 15 261   $null
Update-Row $wsAll 16 79    $null
Update-Row $wsAll 17 1543  $null
Update-Row $wsAll 21 123   $null
Update-Row $wsAll 22 254   $null
Update-Row $wsAll 23 8310  $null
Update-Row $wsAll 25 48    $null
Update-Row $wsAll 28 1285  $null
Update-Row $wsAll 29 40    $null
Update-Row $wsAll 34 6696  $null
Update-Row $wsAll 36 82    $null
Update-Row $wsAll 37 195   $null
Update-Row $wsAll 40 5703  $null
Update-Row $wsAll 44 68    $null

Write-Output "Updates applied."
